# Apply the edits described by the commit:
#  - Remove the "About" sheet note explaining that the model doubles the
#    BLUE Shifts potential (rows 18:19), since the doubling itself is removed.
#  - Remove the "*2" multiplier from the formulas on the PCiCDTdtTDM sheet,
#    so each cell now references the Calcs sheet value directly instead of
#    doubling it.

$wb = $excel.ActiveWorkbook

# --- "About" sheet: delete the now-obsolete explanatory row (and the blank
#     spacer row that followed it), shifting subsequent rows up.
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Rows("18:19").Delete()

# --- "PCiCDTdtTDM" sheet: drop the "*2" doubling from each formula.
$wsPC = $wb.Worksheets.Item("PCiCDTdtTDM")

$wsPC.Range("B2").Formula = "=Calcs!B5"
$wsPC.Range("B3").Formula = "=Calcs!C5"
$wsPC.Range("C3").Formula = "=Calcs!B11"
$wsPC.Range("B4").Formula = "=Calcs!D5"
$wsPC.Range("B5").Formula = "=Calcs!E5"
$wsPC.Range("C5").Formula = "=Calcs!C11"
$wsPC.Range("B6").Formula = "=Calcs!F5"
$wsPC.Range("B7").Formula = "=Calcs!G5"

# B6 and B7 originally carried no explicit cell style (plain "General" format).
# Re-entering their formulas causes the COM layer to auto-inherit the percent
# number format from the cells above (B2:B5), so reset them back to the
# default "Normal" style to match the original (unstyled) formatting.
$wsPC.Range("B6").Style = "Normal"
$wsPC.Range("B7").Style = "Normal"
